# Update the "想去人数" (want-to-go count) figures that changed between
# the two scrape runs. The same three rows are updated on both the
# "展览" sheet and the "全部类型" sheet (they mirror the same events).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 98
    $ws.Range("F5").Value = 2554
    $ws.Range("F6").Value = 237
}
